$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three "VA_CO_*" validation rows (rows 17-19); everything below
# shifts up, which is what the diff shows (rows 20-28 become rows 17-25,
# and the three now-unused strings/styles drop out of sharedStrings/styles).
$ws.Rows("17:19").Delete()

# The row delete doesn't walk into the already-broken (#REF!-based)
# conditional-format formula on E2, so nudge its row window up by 3 to
# match (G$92:G$137 -> G$89:G$134), same as every other reference that
# moved because of the deleted rows.
$fc = $ws.Range("E2").FormatConditions.Item(1)
$fc.Formula1 = "=COUNTIF(G$89:G$134,#REF!)"

# The old last row (now row 25) carried a one-off "applyFill" style that
# nothing else uses; clearing its (already-none) interior pattern drops
# that redundant flag so the cell falls back onto the shared border-only
# style the rest of the column uses.
$ws.Range("C25:D25").Interior.Pattern = -4142

# Leave the selection where the author ended up after the edit.
$ws.Range("B28").Select()
